$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Rows 1-8: direct cell text replacement
$table.Cell(1,1).Range.Text = "53+2=55"
$table.Cell(1,2).Range.Text = "49+42=91"
$table.Cell(1,3).Range.Text = "86-25=61"
$table.Cell(1,4).Range.Text = "40+41=81"
$table.Cell(1,5).Range.Text = "82-49=33"
$table.Cell(2,1).Range.Text = "13+74=87"
$table.Cell(2,2).Range.Text = "29-27=2"
$table.Cell(2,3).Range.Text = "35+49=84"
$table.Cell(2,4).Range.Text = "50-23=27"
$table.Cell(2,5).Range.Text = "44-8=36"
$table.Cell(3,1).Range.Text = "9+57=66"
$table.Cell(3,2).Range.Text = "72-42=30"
$table.Cell(3,3).Range.Text = "25+7=32"
$table.Cell(3,4).Range.Text = "69+11=80"
$table.Cell(3,5).Range.Text = "23-15=8"
$table.Cell(4,1).Range.Text = "15+39=54"
$table.Cell(4,2).Range.Text = "3+5=8"
$table.Cell(4,3).Range.Text = "60-42=18"
$table.Cell(4,4).Range.Text = "85-14=71"
$table.Cell(4,5).Range.Text = "91-29=62"
$table.Cell(5,1).Range.Text = "39+23=62"
$table.Cell(5,2).Range.Text = "24+38=62"
$table.Cell(5,3).Range.Text = "2+53=55"
$table.Cell(5,4).Range.Text = "52+47=99"
$table.Cell(5,5).Range.Text = "18+59=77"
$table.Cell(6,1).Range.Text = "92-27=65"
$table.Cell(6,2).Range.Text = "12-7=5"
$table.Cell(6,3).Range.Text = "45-29=16"
$table.Cell(6,4).Range.Text = "71-39=32"
$table.Cell(6,5).Range.Text = "58-46=12"
$table.Cell(7,1).Range.Text = "10+50=60"
$table.Cell(7,2).Range.Text = "32+62=94"
$table.Cell(7,3).Range.Text = "54-17=37"
$table.Cell(7,4).Range.Text = "24-17=7"
$table.Cell(7,5).Range.Text = "46-17=29"
$table.Cell(8,1).Range.Text = "93-56=37"
$table.Cell(8,2).Range.Text = "37-4=33"
$table.Cell(8,3).Range.Text = "82-9=73"
$table.Cell(8,4).Range.Text = "65-4=61"
$table.Cell(8,5).Range.Text = "37-34=3"

# Row 9: becomes merged content; delete old rows 10-14
$table.Cell(9,1).Range.Text = "5+21=26"
$table.Cell(9,2).Range.Text = "73+21=94"
$table.Cell(9,3).Range.Text = "47+46=93"
$table.Cell(9,4).Range.Text = "45+46=91"
$table.Cell(9,5).Range.Text = "34-25=9"

# Delete rows 10-14 (originally rows 10-14, now collapse since row 9 replaced them)
$table.Rows.Item(10).Delete()
$table.Rows.Item(10).Delete()
$table.Rows.Item(10).Delete()
$table.Rows.Item(10).Delete()
$table.Rows.Item(10).Delete()

# Rows 15-20 (before) are now rows 10-15 after the deletions: direct cell text replacement
$table.Cell(10,1).Range.Text = "67+4=71"
$table.Cell(10,2).Range.Text = "36+33=69"
$table.Cell(10,3).Range.Text = "76-10=66"
$table.Cell(10,4).Range.Text = "14+50=64"
$table.Cell(10,5).Range.Text = "17+51=68"
$table.Cell(11,1).Range.Text = "98-62=36"
$table.Cell(11,2).Range.Text = "74-30=44"
$table.Cell(11,3).Range.Text = "32+7=39"
$table.Cell(11,4).Range.Text = "32+46=78"
$table.Cell(11,5).Range.Text = "66+13=79"
$table.Cell(12,1).Range.Text = "34-14=20"
$table.Cell(12,2).Range.Text = "51-33=18"
$table.Cell(12,3).Range.Text = "15+22=37"
$table.Cell(12,4).Range.Text = "42+2=44"
$table.Cell(12,5).Range.Text = "44+8=52"
$table.Cell(13,1).Range.Text = "39+28=67"
$table.Cell(13,2).Range.Text = "5+1=6"
$table.Cell(13,3).Range.Text = "20+47=67"
$table.Cell(13,4).Range.Text = "31-15=16"
$table.Cell(13,5).Range.Text = "82+1=83"
$table.Cell(14,1).Range.Text = "34+38=72"
$table.Cell(14,2).Range.Text = "72-2=70"
$table.Cell(14,3).Range.Text = "22+67=89"
$table.Cell(14,4).Range.Text = "74-26=48"
$table.Cell(14,5).Range.Text = "92-45=47"
$table.Cell(15,1).Range.Text = "95-19=76"
$table.Cell(15,2).Range.Text = "70-23=47"
$table.Cell(15,3).Range.Text = "92-0=92"
$table.Cell(15,4).Range.Text = "70-12=58"
$table.Cell(15,5).Range.Text = "47-13=34"

# Add 5 new rows at the end
$newRow = $table.Rows.Add()
$newRow = $table.Rows.Add()
$newRow = $table.Rows.Add()
$newRow = $table.Rows.Add()
$newRow = $table.Rows.Add()
$table.Cell(16,1).Range.Text = "61+32=93"
$table.Cell(16,2).Range.Text = "64-43=21"
$table.Cell(16,3).Range.Text = "59-18=41"
$table.Cell(16,4).Range.Text = "60+37=97"
$table.Cell(16,5).Range.Text = "59-44=15"
$table.Cell(17,1).Range.Text = "85+1=86"
$table.Cell(17,2).Range.Text = "42+33=75"
$table.Cell(17,3).Range.Text = "54+13=67"
$table.Cell(17,4).Range.Text = "22+31=53"
$table.Cell(17,5).Range.Text = "60-42=18"
$table.Cell(18,1).Range.Text = "11+17=28"
$table.Cell(18,2).Range.Text = "3+93=96"
$table.Cell(18,3).Range.Text = "1+14=15"
$table.Cell(18,4).Range.Text = "21+59=80"
$table.Cell(18,5).Range.Text = "69-42=27"
$table.Cell(19,1).Range.Text = "29+3=32"
$table.Cell(19,2).Range.Text = "72-2=70"
$table.Cell(19,3).Range.Text = "8+22=30"
$table.Cell(19,4).Range.Text = "15+14=29"
$table.Cell(19,5).Range.Text = "20-11=9"
$table.Cell(20,1).Range.Text = "89-24=65"
$table.Cell(20,2).Range.Text = "92-63=29"
$table.Cell(20,3).Range.Text = "28-23=5"
$table.Cell(20,4).Range.Text = "23+71=94"
$table.Cell(20,5).Range.Text = "29+37=66"
